$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 124, shifting existing rows 124-192 down to 125-193.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new record's data.
$ws.Range("A124").Value = 4
$ws.Range("B124").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C124").Value = "Los Lagos"
$ws.Range("D124").Value = 44518
$ws.Range("E124").Value = 10
$ws.Range("F124").Value = 100112040
$ws.Range("G124").Value = "Cilantro"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 100
$ws.Range("K124").Value = 10000
$ws.Range("L124").Value = 10000
$ws.Range("M124").Value = 10000
$ws.Range("N124").Value = "$/caja 36 atados"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 278
$ws.Range("Q124").Value = 36
$ws.Range("R124").Value = "Hortaliza"
